$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.491361618041992
$ws.Range("B1").Value = 3.787102222442627
$ws.Range("C1").Value = 1.72264575958252
$ws.Range("D1").Value = 1.157914161682129
$ws.Range("E1").Value = 0.749509871006012
